$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# --- Overview sheet: status text for zh-cn / de-de columns ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsOverview.Range("E3").Value = $statusText
$wsOverview.Range("F3").Value = $statusText

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $statusText
$wsZhCn.Range("C3").Value = $statusText

$wsZhCn.Range("I2").Value = "37a9477f-12c8-4502-9d08-c4326715becb.md"
$wsZhCn.Range("I2").Font.Underline = $true
$wsZhCn.Range("I2").Font.Color = 15570276
$wsZhCn.Range("J2").Value = "37a9477f-12c8-4502-9d08-c4326715becb.98d0e4d84b48165b9e71c525f612bf3096d37585.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-08-29 10:22:57"

$wsZhCn.Range("I3").Value = "b4dda690-366d-402c-b78f-e7d58bf5661c.md"
$wsZhCn.Range("I3").Font.Underline = $true
$wsZhCn.Range("I3").Font.Color = 15570276
$wsZhCn.Range("J3").Value = "b4dda690-366d-402c-b78f-e7d58bf5661c.2729becf55cec8b7a6d9be41ec7c264d66af9625.zh-cn.xlf"
$wsZhCn.Range("K3").Value = "2016-08-29 10:22:57"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $statusText
$wsDeDe.Range("C3").Value = $statusText

$wsDeDe.Range("I2").Value = "37a9477f-12c8-4502-9d08-c4326715becb.md"
$wsDeDe.Range("I2").Font.Underline = $true
$wsDeDe.Range("I2").Font.Color = 15570276
$wsDeDe.Range("J2").Value = "37a9477f-12c8-4502-9d08-c4326715becb.98d0e4d84b48165b9e71c525f612bf3096d37585.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-08-29 10:23:11"

$wsDeDe.Range("I3").Value = "b4dda690-366d-402c-b78f-e7d58bf5661c.md"
$wsDeDe.Range("I3").Font.Underline = $true
$wsDeDe.Range("I3").Font.Color = 15570276
$wsDeDe.Range("J3").Value = "b4dda690-366d-402c-b78f-e7d58bf5661c.2729becf55cec8b7a6d9be41ec7c264d66af9625.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-08-29 10:23:11"

Write-Output "done"
